$wb = $excel.ActiveWorkbook

# Add a new "Incidents" worksheet positioned right before the "Accounts" sheet
# (Worksheets.Add(Before) also makes it the active sheet, matching the target
# workbook's new activeTab/tabSelected state).
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("Accounts"))
$ws.Name = "Incidents"

# Re-fetch the Accounts sheet reference AFTER the insert: a reference taken
# beforehand points at the pre-insert sheet position and its later
# `.Select()` would land on the wrong (shifted) sheet.
$acct = $wb.Worksheets.Item("Accounts")

# Populate cells in the same order the shared-string table grew in the
# target file: CaseId, Subject, My Case, Another Case, Company, then the
# values that reuse already-existing shared strings / plain numbers.
$ws.Range("A1").Value = "CaseId"
$ws.Range("B1").Value = "Subject"
$ws.Range("B2").Value = "My Case"
$ws.Range("B3").Value = "Another Case"
$ws.Range("C1").Value = "Company"
$ws.Range("A2").Value = 123
$ws.Range("A3").Value = 456
$ws.Range("C2").Value = "Account A-0001"
$ws.Range("C3").Value = "Account A-0002"

# Column widths for B and C on the new sheet.
$ws.Columns.Item(2).ColumnWidth = 18.291666666666668
$ws.Columns.Item(3).ColumnWidth = 20

# Update the selection on the Accounts sheet (now pushed one tab to the
# right) without leaving it the active sheet.
$acct.Range("B2").Select()

# Re-activate Incidents and park the selection on C6, then make sure it's
# the sheet marked active/tabSelected when saved.
$ws.Activate()
$ws.Range("C6").Select()
